$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = -0.0538
    "E" = -0.527
    "K" = 7.7
    "L" = 0.005268559698939446
    "U" = 2680.1
    "V" = 6.930695629687096
    "W" = 0.005322090129941941
    "X" = 0.2429841320787574
    "Y" = -0.2376620419488154
    "Z" = 0.9664087813264566
    "AB" = 0.1169441081409295
    "AC" = -0.1169441081409295
    "AD" = 1845.6
    "AF" = 1845.6
    "AG" = -834.5
    "AH" = 0.8267705953500875
    "AI" = 0.5012084838280422
    "AJ" = 1.863555158552925
    "AK" = -0.8326681301137497
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
